# Commit message: "Create a separate scenario with definition of missing
# parameter for testing"
#
# This adds a new row (row 6) to the "Scenarios" worksheet describing a new
# test scenario ("TestScenario_missingParam") that references a model
# parameter sheet with a missing parameter ("Global, MissingParam"), mirroring
# the other scenario rows already present on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Scenarios")

# Populate columns in the same left-to-right order used by the other rows
# (A: Scenario_name, B: IndividualId, E: ModelParameterSheets,
#  F: ApplicationProtocol, G: SimulationTime, H: SimulationTimeUnit,
#  L: ModelFile), writing column E before A so that new shared-string
# entries are created in the same order as the target workbook.
$ws.Range("E6").Value = "Global, MissingParam"
$ws.Range("A6").Value = "TestScenario_missingParam"
$ws.Range("B6").Value = "Indiv1"
$ws.Range("F6").Value = "Aciclovir_iv_250mg"
$ws.Range("G6").Value = "0, 24, 60"
$ws.Range("H6").Value = "h"
$ws.Range("L6").Value = "Aciclovir.pkml"

# Move the active selection down to the row below the newly added data.
$ws.Range("E7").Select() | Out-Null
